# Update AVTA course Excel file ("courses" sheet): add the Horticulture /
# Agribusiness Management courses and their combined packages as new rows
# 2-7 (vetCode, cricosCode, department, name, duration, durationDetail,
# tuition, tuitionDetail, location).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Certificate III in Horticulture --------------------------------
$ws.Range("A2").Value = "AHC30716"
$ws.Range("B2").Value = "110597F"
$ws.Range("D2").Value = "CERTIFICATE III IN HORTICULTURE"
$ws.Range("E2").Value = 52
$ws.Range("H2").Value = "44 wks Tuition + 8 wks Break"

# --- Row 3: Certificate IV in Horticulture ----------------------------------
$ws.Range("I2").Value = 10200
$ws.Range("I3").Value = 12200
$ws.Range("J3").Value = "12,000 tuition fee + 200 handling fee"
$ws.Range("J2").Value = "10,000 tuition fee + 200 handling fee"
$ws.Range("M2").Value = "TAS"
$ws.Range("A3").Value = "AHC40416"
$ws.Range("B3").Value = "110598E"
$ws.Range("D3").Value = "CERTIFICATE IV IN HORTICULTURE"
$ws.Range("E3").Value = 52
$ws.Range("H3").Value = "44 wks Tuition + 8 wks Break"
$ws.Range("M3").Value = "TAS"

# --- Row 4: Diploma of Agribusiness Management ------------------------------
$ws.Range("A4").Value = "AHC51422"
$ws.Range("B4").Value = "110774E"
$ws.Range("D4").Value = "DIPLOMA OF AGRIBUSINESS MANAGEMENT"
$ws.Range("E4").Value = 52
$ws.Range("H4").Value = "44 wks Tuition + 8 wks Break"
$ws.Range("I4").Value = 13200
$ws.Range("J4").Value = "13,000 tuition fee + 200 handling fee"
$ws.Range("M4").Value = "TAS"

# --- Row 5: Package - Cert III + Cert IV Horticulture -----------------------
$ws.Range("A5").Value = "AHC30716 / AHC40416"
$ws.Range("B5").Value = "110597F / 110598E"
$ws.Range("D5").Value = "CERTIFICATE III IN HORTICULTURE +`nCERTIFICATE IV IN HORTICULTURE"
$ws.Range("E5").Value = 104
$ws.Range("H5").Value = "88 wks Tuition + 16 wks Break"
$ws.Range("I5").Value = 20200
$ws.Range("J5").Value = "20,000 tuition fee + 200 handling fee"
$ws.Range("M5").Value = "TAS"

# --- Row 6: Package - Cert III Horticulture + Diploma Agribusiness Mgmt ----
$ws.Range("A6").Value = "AHC30716 / AHC51422"
$ws.Range("B6").Value = "110597F / 110774E"
$ws.Range("D6").Value = "CERTIFICATE III IN HORTICULTURE +`nDIPLOMA OF AGRIBUSINESS MANAGEMENT"
$ws.Range("E6").Value = 104
$ws.Range("H6").Value = "88 wks Tuition + 16 wks Break"
$ws.Range("I6").Value = 21200
$ws.Range("J6").Value = "21,000 tuition fee + 200 handling fee"
$ws.Range("M6").Value = "TAS"

# --- Row 7: Package - Cert IV Horticulture + Diploma Agribusiness Mgmt -----
$ws.Range("A7").Value = "AHC40416 / AHC51422"
$ws.Range("B7").Value = "110598E / 110774E"
$ws.Range("D7").Value = "CERTIFICATE IV IN HORTICULTURE +`nDIPLOMA OF AGRIBUSINESS MANAGEMENT"
$ws.Range("E7").Value = 104
$ws.Range("H7").Value = "88 wks Tuition + 16 wks Break"
$ws.Range("I7").Value = 21200
$ws.Range("J7").Value = "21,000 tuition fee + 200 handling fee"
$ws.Range("M7").Value = "TAS"

# --- Department column (filled last) ---------------------------------------
$ws.Range("C2").Value = "HORTICULTURE"
$ws.Range("C3").Value = "HORTICULTURE"
$ws.Range("C4").Value = "MANAGEMENT"
$ws.Range("C5").Value = "PACKAGES"
$ws.Range("C6").Value = "PACKAGES"
$ws.Range("C7").Value = "PACKAGES"

# --- durationDetail column (H) wraps ----------------------------------------
$ws.Range("H2:H7").WrapText = $true

# --- Number format (thousands separator) on tuition column (I) -------------
$ws.Range("I2:I7").NumberFormat = "#,##0"

# --- tuitionDetail column (J) carries the same number format plus wrap -----
$ws.Range("J2:J7").NumberFormat = "#,##0"
$ws.Range("J2:J7").WrapText = $true

# --- vetCode / cricosCode / name wrap for the package rows (5-7) -----------
$ws.Range("A5:B7").WrapText = $true
$ws.Range("D5:D7").WrapText = $true

# --- Row heights --------------------------------------------------------------
$ws.Range("A2:A7").RowHeight = 45

# --- Selection, mirroring the author's final cursor position ---------------
$ws.Range("C10").Select()
